# MCH150.xlsx: add the MK AND THE FUTURE record as row 2 of Sheet1.
# (commit message: "Updated MCH102 to MCH251" — read literally this is
#  about correcting/adding the MCH150 collection-level metadata row.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 values -----------------------------------------------------
$ws.Range("A2").Value = "MCH150"
$ws.Range("C2").Value = "MK AND THE FUTURE, BORN OF THE PEOPLE, JOURNAL OF THE MK"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"
# D2 and H2 carry the row's formatting but stay empty, like their
# neighbours (no identifier alternative / file_path for this record).

# --- Formatting: identifier / title / levelOfDescription / extent /
#     notes cells all pick up the small 10pt Calibri "data row" font
#     (automatic/theme text colour) used throughout the rest of the
#     sheet's body rows. ---------------------------------------------
$cells = "A2", "C2", "D2", "E2", "F2", "G2", "H2"
foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 10
    $rng.Font.ThemeColor = 1
}
$ws.Range("F2").WrapText = $false

# Keep the explicit (already-default) row height, matching how Excel
# pins down ht/customHeight once a sheet is touched and resaved.
$ws.Rows.Item(2).RowHeight = 15.75

# --- View state: select the new row and re-freeze the header row -----
$ws.Range("A2:I2").Select()
$excel.ActiveWindow.FreezePanes = $true
